$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Add the missing shared-string cell in row 3 (same value already used in C2)
$ws1.Range("C3").Value = "accepts=text/html; charset=utf-8"

# Update zoom on both sheets and change the active selection on sheet 1
$ws1.Activate()
$excel.ActiveWindow.Zoom = 160
$ws1.Range("H10").Select() | Out-Null

$ws2.Activate()
$excel.ActiveWindow.Zoom = 160
$ws2.Range("J1").Select() | Out-Null

$ws1.Activate()
